$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $value) {
    $cell = $ws.Range($cellref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "62.623.06"
$ws.Range("E2").Value = "  -0.61%  "

Set-TextValue "D3" "2.575.33"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "580.68"
$ws.Range("E5").Value = "  -0.30%  "

Set-TextValue "D6" "144.79"
$ws.Range("E6").Value = "  -1.32%  "

$ws.Range("E7").Value = "  -0.03%  "

Set-TextValue "D8" "0.591"
$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("E9").Value = "  +0.31%  "

Set-TextValue "D10" "5.56"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("E12").Value = "  -0.90%  "

Set-TextValue "D13" "26.91"
$ws.Range("E13").Value = "  -2.10%  "

Set-TextValue "D14" "3.035.96"
$ws.Range("E14").Value = "  +0.99%  "

Set-TextValue "D15" "62.596.42"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("E16").Value = "  +0.54%  "

Set-TextValue "D17" "2.567.43"
$ws.Range("E17").Value = "  +0.94%  "

Set-TextValue "D18" "11.17"
$ws.Range("E18").Value = "  -1.41%  "

Set-TextValue "D19" "337.68"
$ws.Range("E19").Value = "  -0.23%  "

Set-TextValue "D20" "4.34"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("E22").Value = "  +0.10%  "

Set-TextValue "D23" "66.92"
$ws.Range("E23").Value = "  +1.99%  "

Set-TextValue "D24" "2.699.21"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  -1.21%  "

Set-TextValue "D29" "7.83"
$ws.Range("E29").Value = "  +2.10%  "

Set-TextValue "D30" "8.19"
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D32" "461.37"
$ws.Range("E32").Value = "  +10.47%  "

$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D33" "0.0₃0807"
$ws.Range("E33").Value = "  -0.83%  "

Set-TextValue "D34" "177.02"
$ws.Range("E34").Value = "  -0.59%  "

Set-TextValue "D35" "1.59"
$ws.Range("E35").Value = "  +3.48%  "

$ws.Range("E36").Value = "  +0.01%  "

Set-TextValue "D37" "0.400"
$ws.Range("E37").Value = "  -0.21%  "

Set-TextValue "D38" "18.87"
$ws.Range("E38").Value = "  -1.17%  "

Set-TextValue "D39" "4.46"
$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  -3.41%  "

Set-TextValue "D42" "157.07"
$ws.Range("E42").Value = "  +4.22%  "

$ws.Range("E43").Value = "  -1.69%  "

Set-TextValue "D44" "21.03"
$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("E45").Value = "  +3.89%  "

Set-TextValue "D46" "0.0533"
$ws.Range("E46").Value = "  -0.86%  "

$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("E48").Value = "  -1.59%  "

Set-TextValue "D49" "18.06"
$ws.Range("E49").Value = "  -1.29%  "

$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("E51").Value = "  -0.47%  "
